$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9, shifting existing rows 9-12 down to 10-13.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value = 12
$ws.Cells.Item(9, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44467
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 35
$ws.Cells.Item(9, 11).Value = 12000
$ws.Cells.Item(9, 12).Value = 12000
$ws.Cells.Item(9, 13).Value = 12000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 480
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Ensure the date column keeps the same numeric/date style as the rest of column D.
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
